# Applies two changes captured in the target commit:
#
#  1. The cash-flow table on slide 16 is switched from the deck's custom
#     "Table_0" table style to the built-in table style
#     {A1877012-033E-459F-98BF-2AD2F094BFE9}.
#
#  2. The presentation's theme (ppt/theme/theme1.xml, the theme used by the
#     slide master / design "Integral") is recoloured to the stock Office
#     theme palette ("Office Theme" / clrScheme "Office") -- i.e. the
#     Integral and Office colour schemes trade places between the two theme
#     parts in the package. The PowerPoint object model only exposes the
#     slide-master's theme colours (ThemeColorScheme.Colors(i).RGB), so that
#     is the side of the swap performed here.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$slide = $p.Slides.Item(16)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A1877012-033E-459F-98BF-2AD2F094BFE9}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
# VBA RGB() packs a hex colour "RRGGBB" as (B*65536 + G*256 + R).
function HexToVbaRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Target palette = the stock "Office Theme" colour scheme.
$officeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

$themeColorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Length; $i++) {
    $themeColorScheme.Colors($i).RGB = HexToVbaRgb $officeColors[$i - 1]
}
